# 4.0.3 model and data
#
# The "Boolean" sheet lists InputData CSV files whose cells must contain a
# value constrained to a specific set ("boolean"-style list). The single
# "trans/BVTQaZ/BVTQaZ.csv" and "trans/VTQaZ/VTQaZ.csv" aggregate files were
# split into six per-vehicle-type files apiece (LDVs/HDVs/aircraft/rail/
# ships/motorbikes), and a handful of blank rows were left below the list.
# Sheet selections / the active tab are also refreshed to match the
# authored workbook state.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Boolean" sheet: split the BVTQaZ and VTQaZ rows into 6 rows each.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Boolean")

# --- expand "trans/BVTQaZ/BVTQaZ.csv" (row 17) into 6 rows ---
$ws3.Rows("18:22").Insert()
$ws3.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$ws3.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$ws3.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$ws3.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$ws3.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$ws3.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# after the insert above, "trans/VTQaZ/VTQaZ.csv" (originally row 21)
# now lives at row 26 -- expand it into 6 rows too.
$ws3.Rows("27:31").Insert()
$ws3.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$ws3.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$ws3.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$ws3.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$ws3.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$ws3.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# a handful of blank, body-formatted rows trail the list (rows 33-38)
$ws3.Rows("33:38").Font.Name = "Calibri"
$ws3.Rows("33:38").Font.Size = 11

# scroll position / selection left on this sheet
$ws3.Range("A32").Select() | Out-Null

# ---------------------------------------------------------------------
# Restore sheet selections / active tab as saved in the workbook.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Integer")
$ws2.Range("A13").Select() | Out-Null

$ws1 = $wb.Worksheets.Item("About")
$ws1.Activate() | Out-Null
